$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- Row 17: begin group "cholera_cases" ------------------------------
$ws.Range("A17").Value = "begin group"
$ws.Range("B17").Value = "cholera_cases"
$ws.Range("C17").Value = "Cholera Cases"

# D17 and F17 are blank cells that merely carry the default style (s=1),
# matching the pattern already used on row 7 (G7/H7). Copy format only
# from an existing default-styled blank cell so no new style is created.
$ws.Range("G7").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("F17").PasteSpecial(-4122)

# --- Row 18: num_screened (uses the shaded/emphasised style s=2) ------
$ws.Range("A18").Value = "integer"
$ws.Range("B18").Value = "num_screened"

$ws.Range("C15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("F18").PasteSpecial(-4122)

$ws.Range("C18").Value = "Number Of Persons Screened For Cholera"
$ws.Range("D18").Value = "yes"
$ws.Range("F18").Value = "numbers"

# --- Row 19: num_referred ----------------------------------------------
$ws.Range("A19").Value = "integer"
$ws.Range("B19").Value = "num_referred"
$ws.Range("C19").Value = "Number Of Presumptive Cholera Persons Referred For Diagnosis"
$ws.Range("D19").Value = "yes"
$ws.Range("F19").Value = "numbers"

# --- Row 20: num_referred_reached (has relevant expr in column E) -----
$ws.Range("A20").Value = "integer"
$ws.Range("B20").Value = "num_referred_reached"
$ws.Range("C20").Value = "Number Of Referred Persons Who Reached Health Facility"
$ws.Range("D20").Value = "yes"
$ws.Range("E20").Value = '${num_referred} >0'
$ws.Range("F20").Value = "numbers"

# --- Row 21: num_confirmed_cases ---------------------------------------
$ws.Range("A21").Value = "integer"
$ws.Range("B21").Value = "num_confirmed_cases"
$ws.Range("C21").Value = "Number Of Confirmed Cholera Cases At Health Facility"
$ws.Range("D21").Value = "yes"
$ws.Range("F21").Value = "numbers"

# --- Row 22: num_deaths -------------------------------------------------
$ws.Range("A22").Value = "integer"
$ws.Range("B22").Value = "num_deaths"
$ws.Range("C22").Value = "Number Of Deaths Due To Cholera In The Month"
$ws.Range("D22").Value = "yes"
$ws.Range("F22").Value = "numbers"

# --- Row 23: end group ---------------------------------------------------
$ws.Range("A23").Value = "end group"

# --- Column E needs an explicit width now that it carries content -----
$ws.Columns.Item(5).ColumnWidth = 21.2

Write-Output "cholera cases group added"
